# Update the old ATACseq
# This script applies the v2.0.0 ATACseq schema changes:
#  - dataset_type: remove "RNAseq (Visium)", "GeoMx", "RNAseq (GeoMx)";
#                  add "GeoMx (NGS)" (after CODEX) and "GeoMx (nCounter)" (after 10X Multiome)
#  - library_concentration_unit: add "nM" (http://purl.obolibrary.org/obo/UO_0000065)
#  - umi_offset: add "0" as a new first option
#  - .metadata: bump pav:createdOn timestamp
#  - update the 3 affected data validation list ranges on the ATACseq sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. dataset_type sheet: remove 3 rows, add 2 rows (35 -> 34 options)
# ---------------------------------------------------------------------------
$wsDataset = $wb.Worksheets.Item("dataset_type")

# Delete from bottom to top so earlier row numbers stay valid:
#   row 24 = RNAseq (GeoMx)
#   row 20 = GeoMx
#   row 15 = RNAseq (Visium)
$wsDataset.Rows.Item(24).Delete()
$wsDataset.Rows.Item(20).Delete()
$wsDataset.Rows.Item(15).Delete()

# After the deletes the list (1-based) is:
#   ... 13 DBiT, 14 SIMS, 15 Cell DIVE, 16 CODEX, 17 CyCIF, 18 Light Sheet,
#   19 RNAseq (bulk), ... 27 LC-MS, 28 10X Multiome, 29 PhenoCycler, ...
# Insert "GeoMx (nCounter)" right after "10X Multiome" (new row 28) -> becomes row 29
$wsDataset.Rows.Item(29).Insert()
$wsDataset.Cells.Item(29,1).Value = "GeoMx (nCounter)"
$wsDataset.Cells.Item(29,2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000301"

# Insert "GeoMx (NGS)" right after "CODEX" (row 16) -> becomes row 17
$wsDataset.Rows.Item(17).Insert()
$wsDataset.Cells.Item(17,1).Value = "GeoMx (NGS)"
$wsDataset.Cells.Item(17,2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000300"

# ---------------------------------------------------------------------------
# 2. library_concentration_unit sheet: add "nM" as a 2nd option
# ---------------------------------------------------------------------------
$wsLibConc = $wb.Worksheets.Item("library_concentration_unit")
$wsLibConc.Cells.Item(2,1).Value = "nM"
$wsLibConc.Cells.Item(2,2).Value = "http://purl.obolibrary.org/obo/UO_0000065"

# ---------------------------------------------------------------------------
# 3. umi_offset sheet: add "0" as the new first option (push others down)
#    Copy the existing text cell "0" from barcode_offset!A1 so that it is
#    stored as a real text value (matching the shared string "0") instead
#    of being auto-converted to a number.
# ---------------------------------------------------------------------------
$wsUmiOffset = $wb.Worksheets.Item("umi_offset")
$wsBarcodeOffset = $wb.Worksheets.Item("barcode_offset")

$wsUmiOffset.Rows.Item(1).Insert()
$wsBarcodeOffset.Range("A1").Copy()
$wsUmiOffset.Range("A1").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. .metadata sheet: bump the pav:createdOn timestamp
# ---------------------------------------------------------------------------
$wsMetadata = $wb.Worksheets.Item(".metadata")
$wsMetadata.Cells.Item(2,3).Value = "2023-11-22T10:19:45-08:00"

# ---------------------------------------------------------------------------
# 5. Update the data validation list ranges that shifted because of the
#    row count changes above.
# ---------------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("ATACseq")

$dvDataset = $wsMain.Range("D2:D1001").Validation
$dvDataset.Modify(3, 1, 1, "='dataset_type'!`$A`$1:`$A`$34")

$dvUmiOffset = $wsMain.Range("R2:R1001").Validation
$dvUmiOffset.Modify(3, 1, 1, "='umi_offset'!`$A`$1:`$A`$3")

$dvLibConc = $wsMain.Range("AD2:AD1001").Validation
$dvLibConc.Modify(3, 1, 1, "='library_concentration_unit'!`$A`$1:`$A`$2")
